# Add team win/loss/tie record columns to the player data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold, centered, bordered) from the last
# existing header cell (AC1) onto the three new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# New header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-46) gets the same team record for the 2002 season.
$ws.Range("AD2:AD46").Value = 97
$ws.Range("AE2:AE46").Value = 65
$ws.Range("AF2:AF46").Value = 0

Write-Output "Added Wins/Losses/Ties columns (AD:AF) for rows 1-46"
